# Set up the print/view configuration on the "Eng_dict" sheet entry:
#  - Repair the broken Print_Area defined name (#REF!) to cover the
#    dictionary table, B3:K28 plus header/margin -> A1:K29
#  - Select that same range so it becomes the sheet's live selection
#  - Turn on printed row/column headings and gridlines
#  - Re-scale the print layout from 73% to 59%

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eng_dict")

$ws.Activate()

# Print area: Eng_dict!$A$1:$K$29 (was Eng_dict!#REF!)
$ws.PageSetup.PrintArea = '$A$1:$K$29'

# Select A1:K29 on the sheet (was B3:K28)
$ws.Range("A1:K29").Select()

# Print row/column headings and gridlines
$ws.PageSetup.PrintHeadings = $true
$ws.PageSetup.PrintGridlines = $true

# Print scale: 59% (was 73%), keep "fit to 2 pages tall"
$ws.PageSetup.Zoom = 59
$ws.PageSetup.FitToPagesTall = 2
